# Apply the "Add files via upload" edit: add a new sample row (SS-10)
# and fill in a few missing values in existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank E column values for rows 6, 9, 10
$ws.Range("E6").Value = 1.7
$ws.Range("E9").Value = 1.8
$ws.Range("E10").Value = 1.8

# Fill in previously-blank N column values for rows 11, 12
$ws.Range("N11").Value = 240
$ws.Range("N12").Value = 240

# Add new row 13 data (SS-10 sample)
$ws.Range("A13").Value = "SS-10"
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 9.1
$ws.Range("E13").Value = 1.9
$ws.Range("N13").Value = 300

# Update selection to match the post-edit state
$ws.Range("M15").Select()
